$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Wong3")
$ws.Name = "Euclid"

$ws.Range("C2").Value = 228
$ws.Range("D2").Value = 92.30769230769231
$ws.Range("C3").Value = 14
$ws.Range("D3").Value = 5.668016194331984
$ws.Range("C4").Value = 44
$ws.Range("D4").Value = 17.81376518218623
$ws.Range("C5").Value = 16
$ws.Range("D5").Value = 6.477732793522267
$ws.Range("C6").Value = 228
$ws.Range("D6").Value = 92.30769230769231
$ws.Range("C7").Value = 42
$ws.Range("D7").Value = 17.00404858299595
$ws.Range("C8").Value = 4
$ws.Range("D8").Value = 1.619433198380567
$ws.Range("C9").Value = 6
$ws.Range("D9").Value = 2.42914979757085
$ws.Range("C10").Value = 227
$ws.Range("D10").Value = 91.90283400809717
$ws.Range("C11").Value = 110
$ws.Range("D11").Value = 44.53441295546558
$ws.Range("C12").Value = 223
$ws.Range("D12").Value = 90.2834008097166
$ws.Range("C13").Value = 1
$ws.Range("D13").Value = 0.4048582995951417
$ws.Range("C14").Value = 2
$ws.Range("D14").Value = 0.8097165991902834
$ws.Range("C15").Value = 75
$ws.Range("D15").Value = 30.36437246963563
$ws.Range("C16").Value = 12
$ws.Range("D16").Value = 4.8582995951417
$ws.Range("C17").Value = 17
$ws.Range("D17").Value = 6.882591093117409
$ws.Range("C18").Value = 2
$ws.Range("D18").Value = 0.8097165991902834
$ws.Range("C19").Value = 92
$ws.Range("D19").Value = 37.24696356275304
$ws.Range("C20").Value = 5
$ws.Range("D20").Value = 2.024291497975709
$ws.Range("C21").Value = 87
$ws.Range("D21").Value = 35.22267206477733
$ws.Range("C22").Value = 3
$ws.Range("D22").Value = 1.214574898785425
$ws.Range("C23").Value = 228
$ws.Range("D23").Value = 92.30769230769231
$ws.Range("C24").Value = 27
$ws.Range("D24").Value = 10.93117408906883
$ws.Range("C25").Value = 2
$ws.Range("D25").Value = 0.8097165991902834
$ws.Range("C26").Value = 24
$ws.Range("D26").Value = 9.7165991902834
$ws.Range("C27").Value = 2
$ws.Range("D27").Value = 0.8097165991902834
$ws.Range("C28").Value = 228
$ws.Range("D28").Value = 92.30769230769231
$ws.Range("C29").Value = 12
$ws.Range("D29").Value = 4.8582995951417
$ws.Range("C30").Value = 18
$ws.Range("D30").Value = 7.28744939271255
$ws.Range("C31").Value = 2
$ws.Range("D31").Value = 0.8097165991902834
$ws.Range("C32").Value = 97
$ws.Range("D32").Value = 39.27125506072874
$ws.Range("C33").Value = 40
$ws.Range("D33").Value = 16.19433198380567
$ws.Range("C34").Value = 228
$ws.Range("D34").Value = 92.30769230769231
$ws.Range("C35").Value = 2
$ws.Range("D35").Value = 0.8097165991902834
$ws.Range("C36").Value = 4
$ws.Range("D36").Value = 1.619433198380567
$ws.Range("C37").Value = 42
$ws.Range("D37").Value = 17.00404858299595
$ws.Range("C38").Value = 97
$ws.Range("D38").Value = 39.27125506072874
$ws.Range("C39").Value = 14
$ws.Range("D39").Value = 5.668016194331984
$ws.Range("C40").Value = 3
$ws.Range("D40").Value = 1.214574898785425
$ws.Range("C41").Value = 20
$ws.Range("D41").Value = 8.097165991902834
$ws.Range("C42").Value = 97
$ws.Range("D42").Value = 39.27125506072874
$ws.Range("C43").Value = 12
$ws.Range("D43").Value = 4.8582995951417
$ws.Range("C44").Value = 3
$ws.Range("D44").Value = 1.214574898785425
$ws.Range("C45").Value = 54
$ws.Range("D45").Value = 21.86234817813765
$ws.Range("C46").Value = 228
$ws.Range("D46").Value = 92.30769230769231
$ws.Range("C47").Value = 58
$ws.Range("D47").Value = 23.48178137651822
$ws.Range("C48").Value = 228
$ws.Range("D48").Value = 92.30769230769231
$ws.Range("C49").Value = 97
$ws.Range("D49").Value = 39.27125506072874
$ws.Range("C50").Value = 12
$ws.Range("D50").Value = 4.8582995951417
$ws.Range("C51").Value = 46
$ws.Range("D51").Value = 18.62348178137652
$ws.Range("C52").Value = 97
$ws.Range("D52").Value = 39.27125506072874
$ws.Range("C53").Value = 58
$ws.Range("D53").Value = 23.48178137651822
$ws.Range("C54").Value = 227
$ws.Range("D54").Value = 91.90283400809717
$ws.Range("C55").Value = 53
$ws.Range("D55").Value = 21.45748987854251
$ws.Range("C56").Value = 41
$ws.Range("D56").Value = 16.59919028340081
$ws.Range("C57").Value = 54
$ws.Range("D57").Value = 21.86234817813765
$ws.Range("C58").Value = 12
$ws.Range("D58").Value = 4.8582995951417
$ws.Range("C59").Value = 62
$ws.Range("D59").Value = 25.10121457489878
$ws.Range("C60").Value = 227
$ws.Range("D60").Value = 91.90283400809717
$ws.Range("C61").Value = 12
$ws.Range("D61").Value = 4.8582995951417
$ws.Range("C62").Value = 2
$ws.Range("D62").Value = 0.8097165991902834
$ws.Range("C63").Value = 27
$ws.Range("D63").Value = 10.93117408906883
$ws.Range("C64").Value = 227
$ws.Range("D64").Value = 91.90283400809717
$ws.Range("C65").Value = 8
$ws.Range("D65").Value = 3.238866396761134
$ws.Range("C66").Value = 227
$ws.Range("D66").Value = 91.90283400809717
$ws.Range("C67").Value = 71
$ws.Range("D67").Value = 28.74493927125506
$ws.Range("C68").Value = 58
$ws.Range("D68").Value = 23.48178137651822
$ws.Range("C69").Value = 14
$ws.Range("D69").Value = 5.668016194331984
